$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.307.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.87%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.128.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.60%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.35%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "611.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "

# Row 7
$ws.Range("E7").Value = "  +1.86%  "

# Row 8
$ws.Range("E8").Value = "  -1.09%  "

# Row 9
$ws.Range("E9").Value = "  +0.04%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.127.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.61%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.781"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.46%  "

# Row 12
$ws.Range("E12").Value = "  -0.16%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.866.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.70%  "

# Row 14
$ws.Range("E14").Value = "  -1.58%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.18%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.96%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.711.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.77%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.127.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.33%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "530.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +20.66%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.48%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.41%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.84%  "

# Row 23
$ws.Range("E23").Value = "  -5.18%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.95%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.86%  "

# Row 26
$ws.Range("E26").Value = "  -2.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.61%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.297.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.71%  "

# Row 29
$ws.Range("E29").Value = "  -0.19%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.238"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.94%  "

# Row 31
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.175"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.99%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.125"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.24%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.56%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.74%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.14%  "

# Row 36
$ws.Range("E36").Value = "  -6.23%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.82%  "

# Row 38
$ws.Range("E38").Value = "  -0.69%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "484.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.89%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.37%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.441"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.04%  "

# Row 42
$ws.Range("E42").Value = "  -4.64%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.78%  "

# Row 44
$ws.Range("E44").Value = "  -0.02%  "

# Row 45
$ws.Range("E45").Value = "  -4.99%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "160.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.04%  "

# Row 47
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.86%  "

# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.701"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.74%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.12%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.29%  "

# Row 51
$ws.Range("E51").Value = "  +0.09%  "
